$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The five data blocks (rows 12-16, 17-21, 22-26, 27-31, 32-36) each hold a
# full record (columns A:D) for one "version" value. The edit rotates the
# block contents down by one block (the last block's content wraps around
# to become the first block's content), while columns B and D stay the
# same because their per-offset pattern repeats identically in every block.

# Capture the original column A (version label) and column C (value) for
# each block before overwriting, keyed by block start row.
$blockStarts = @(12, 17, 22, 27, 32)

$origA = @{}
$origC = @{}
foreach ($start in $blockStarts) {
    $a = @()
    $c = @()
    for ($i = 0; $i -lt 5; $i++) {
        $r = $start + $i
        $a += $ws.Cells.Item($r, 1).Value2
        $c += $ws.Cells.Item($r, 3).Value2
    }
    $origA[$start] = $a
    $origC[$start] = $c
}

# New block order: block at index k receives the content that was
# originally in the block before it (wrapping), i.e. block4 -> block0,
# block0 -> block1, block1 -> block2, block2 -> block3, block3 -> block4.
$srcForDest = @{
    12 = 32
    17 = 12
    22 = 17
    27 = 22
    32 = 27
}

foreach ($destStart in $blockStarts) {
    $srcStart = $srcForDest[$destStart]
    $aVals = $origA[$srcStart]
    $cVals = $origC[$srcStart]
    for ($i = 0; $i -lt 5; $i++) {
        $r = $destStart + $i
        $ws.Cells.Item($r, 1).Value2 = $aVals[$i]
        $ws.Cells.Item($r, 3).Value2 = $cVals[$i]
    }
}

# Update the active selection on the sheet.
$ws.Range("E8").Select()
